$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3, 0.02214582616087132, 0.01392707507159257, 0.1760912821103572, 0.06030058583506512),
    @(4, 0.02062567101747944, 0.01295000086963814, 0.1703176609672555, 0.05618998089931112),
    @(6, 0.02040793302508948, 0.01343025013593323, 0.1684352693398844, 0.05924301678028802),
    @(5, 0.01836442606289615, 0.01041808638625798, 0.1618114638400777, 0.05021768464686466),
    @(7, 0.01964345860804612, 0.01191061316394075, 0.1666940905568095, 0.05471171046745366),
    @(8, 0.01801037019193109, 0.0108761227886058, 0.159066785795196, 0.0513337030217291)
)

$startRow = 4
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le 5; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}
